$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103; this shifts the existing rows 103:186
# down to 104:187 (carrying their values/styles with them), and extends
# the sheet dimension to A1:R187 automatically.
$ws.Rows(103).Insert()

# Populate the newly inserted row 103 with the new price record.
$ws.Range("A103").Value = 3
$ws.Range("B103").Value = "Femacal de La Calera"
$ws.Range("C103").Value = "Coquimbo"
$ws.Range("D103").Value = 44447
$ws.Range("E103").Value = 5
$ws.Range("F103").Value = 100112009
$ws.Range("G103").Value = "Acelga"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 230
$ws.Range("K103").Value = 2000
$ws.Range("L103").Value = 2300
$ws.Range("M103").Value = 2143
$ws.Range("N103").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O103").Value = "Provincia de Quillota"
$ws.Range("P103").Value = 357
$ws.Range("Q103").Value = 6
$ws.Range("R103").Value = "Hortaliza"
